$d = $word.ActiveDocument

# --- 1. Body sectPr: restart footnote numbering each section -----------
# <w:sectPr/>  ->  <w:sectPr><w:footnotePr><w:numRestart w:val="eachSect"/></w:footnotePr></w:sectPr>
# wdRestartEachSection = 1
$d.Footnotes.NumberingRule = 1

# --- 2. Title / TitleChar styles: drop the expanded/kerned look --------
# (remove <w:spacing w:val="-10"/> and <w:kern w:val="28"/> from rPr)
$titleStyle = $d.Styles("Title")
$titleStyle.Font.Spacing = 0
$titleStyle.Font.Kerning = 0

$titleCharStyle = $d.Styles("TitleChar")
$titleCharStyle.Font.Spacing = 0
$titleCharStyle.Font.Kerning = 0

# --- 3. Author / Date styles: base them on Title, drop the explicit ----
#        center alignment (inherited from Title) and give them their
#        own smaller run size (12pt / 12pt cs).
$authorStyle = $d.Styles("Author")
$authorStyle.BaseStyle = $titleStyle
$authorStyle.ParagraphFormat.Alignment = 1
$authorStyle.Font.Size = 12
$authorStyle.Font.SizeBi = 12

$dateStyle = $d.Styles("Date")
$dateStyle.BaseStyle = $titleStyle
$dateStyle.ParagraphFormat.Alignment = 1
$dateStyle.Font.Size = 12
$dateStyle.Font.SizeBi = 12
